# Apply Sat May  4 09:44:52 UTC 2024 cryptos-list refresh (GitHub Actions bot commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.247.66'
$ws.Range("E2").Value = '  +6.30%  '
$ws.Range("D3").Value = '3.111.59'
$ws.Range("E3").Value = '  +4.24%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'584.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.27%  '
$ws.Range("D6").Value = "'144.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.41%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.102.27'
$ws.Range("E8").Value = '  +4.25%  '
$ws.Range("E9").Value = '  +1.87%  '
$ws.Range("E10").Value = '  +13.30%  '
$ws.Range("D11").Value = "'5.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.02%  '
$ws.Range("D12").Value = "'0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.64%  '
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.58%  '
$ws.Range("D14").Value = "'35.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.30%  '
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("D16").Value = '3.626.03'
$ws.Range("E16").Value = '  +4.27%  '
$ws.Range("D17").Value = "'7.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("D18").Value = '63.159.29'
$ws.Range("E18").Value = '  +6.18%  '
$ws.Range("D19").Value = '3.110.76'
$ws.Range("E19").Value = '  +4.32%  '
$ws.Range("D20").Value = "'466.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.81%  '
$ws.Range("E21").Value = '  +3.70%  '
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("D23").Value = "'7.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.01%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = "'82.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.47%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = "'8.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.90%  '
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("B29").Value = 'FirstDigitalUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'2.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.66%  '
$ws.Range("D31").Value = "'6.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.53%  '
$ws.Range("D32").Value = "'26.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.49%  '
$ws.Range("E33").Value = '  +4.22%  '
$ws.Range("D34").Value = '0.0₃0864'
$ws.Range("E34").Value = '  +11.24%  '
$ws.Range("D35").Value = "'2.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +16.04%  '
$ws.Range("E36").Value = '  +7.06%  '
$ws.Range("E37").Value = '  +2.69%  '
$ws.Range("D38").Value = "'3.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +19.15%  '
$ws.Range("D39").Value = "'50.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.85%  '
$ws.Range("D40").Value = "'439.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.06%  '
$ws.Range("D41").Value = "'8.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.51%  '
$ws.Range("D42").Value = '2.918.17'
$ws.Range("E42").Value = '  +6.83%  '
$ws.Range("D43").Value = "'0.0369"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.01%  '
$ws.Range("D44").Value = "'0.279"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.19%  '
$ws.Range("E45").Value = '  +4.09%  '
$ws.Range("E46").Value = '  +7.86%  '
$ws.Range("D48").Value = "'34.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.60%  '
$ws.Range("D49").Value = "'123.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.74%  '
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("D51").Value = "'24.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.69%  '
